$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Update the "Section" header row (was "Header") - content only, same row
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = "Section"

# ---------------------------------------------------------------------------
# 2) Update row 4 (was "NA"/8/long-combined-text) to become the first
#    "Carousel" row with the first bullet "1)Surge Protect"
# ---------------------------------------------------------------------------
$ws.Range("A4").Value = "Carousel"
$ws.Range("B4").Value = 7
$ws.Range("C4").Value = "1)Surge Protect"

# ---------------------------------------------------------------------------
# 3) Insert 6 new rows right after row 4 (rows 5-10) to hold the remaining
#    bullet points for the "Carousel" / "Home(Dashboard)" section.
#    This pushes the old row 5 ("Discover" merged header) down to row 11.
# ---------------------------------------------------------------------------
$ws.Range("A5:A10").EntireRow.Insert()
$ws.Range("A4:C4").Copy()
$ws.Range("A5:C10").PasteSpecial(-4122)
$ws.Range("A5:B10").ClearContents()

$ws.Range("C5").Value = "2)AC and Heat Protect"
$ws.Range("C6").Value = "3)DoD"
$ws.Range("C7").Value = "4)Average Billing"
$ws.Range("C8").Value = "5)Electrical line and Surge Protect"
$ws.Range("C9").Value = "6)Make it Solar"
$ws.Range("C10").Value = "7)Reliant Ecoshare"

# ---------------------------------------------------------------------------
# 4) Row 11 is now "Discover" (merged) - unchanged content, just moved.
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# 5) Replace the old single-line "Protection Plans" / "Renewable Products" /
#    "Home Automation and Security" / "Backup Power" rows (old rows 6-9,
#    now shifted to rows 12-15) with the new expanded layout:
#      Renewable Products (3 promos)   -> rows 12-14
#      Protection Plans   (5 promos)   -> rows 15-19
#      Backup Power       (1 promo)    -> row 20
#      Home Automation... (1 promo)    -> row 21
#    First, insert enough additional rows to fit everything.
#    Currently rows 12-15 hold the 4 old category rows (1 row each = 4 rows).
#    We need 3 + 5 + 1 + 1 = 10 rows total, so insert 6 more rows after the
#    current row 15 (old row 9, "Backup Power").
# ---------------------------------------------------------------------------
$ws.Range("A16:A21").EntireRow.Insert()
$ws.Range("A12:C12").Copy()
$ws.Range("A12:C21").PasteSpecial(-4122)
$ws.Range("A12:B21").ClearContents()
$ws.Range("C12:C21").ClearContents()

# Renewable Products (rows 12-14)
$ws.Range("A12").Value = "Renewable Products"
$ws.Range("B12").Value = 3
$ws.Range("C12").Value = "1)Make it Solar"
$ws.Range("C13").Value = "2)Reliant EcoShare"
$ws.Range("C14").Value = "3)DoD"

# Protection Plans (rows 15-19)
$ws.Range("A15").Value = "Protection Plans"
$ws.Range("B15").Value = 5
$ws.Range("C15").Value = "1)Surge Protect"
$ws.Range("C16").Value = "2)AC and Heat Protect"
$ws.Range("C17").Value = "3)Electric Protect"
$ws.Range("C18").Value = "4)Plumbing Protect"
$ws.Range("C19").Value = "5)Electric Line and Surge Protect"

# Backup Power (row 20)
$ws.Range("A20").Value = "Backup Power"
$ws.Range("B20").Value = 1
$ws.Range("C20").Value = "1)Whole-home Generators"

# Home Automation and Security (row 21)
$ws.Range("A21").Value = "Home Automation and Security"
$ws.Range("B21").Value = 1
$ws.Range("C21").Value = "1)Vivint"

# ---------------------------------------------------------------------------
# 6) Column C width - narrower now that entries are shorter.
# ---------------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 28.25

# ---------------------------------------------------------------------------
# 7) Update the saved selection to match the authored file.
# ---------------------------------------------------------------------------
$ws.Range("N15").Select()
